$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "maa://40192 (98.15), maa://36987 (96.08), maa://39849 (88.89)"
$ws.Range("X4").Value = "**maa://32495 (48.7), ***maa://31785 (22.22), maa://43217 (90.48), ***maa://36683 (28.26)"
$ws.Range("AF6").Value = "*maa://33152 (64.15), ***maa://22770 (26.09)"
$ws.Range("A8").Value = "更新日期：2025.03.01 13:17:29"
$ws.Range("P8").Value = "maa://32931 (83.9), *maa://21916 (62.12), maa://23252 (91.18), maa://37496 (96.97), **maa://22759 (45.45)"
$ws.Range("D10").Value = "***maa://25695 (18.62), ***maa://39951 (14.81), ***maa://34206 (20.0), ***maa://39243 (25.0), *maa://45271 (57.89)"
$ws.Range("AF10").Value = "*maa://25021 (54.17), *maa://22733 (61.11), **maa://22761 (50.0)"
$ws.Range("X12").Value = "maa://22753 (91.01), *maa://21485 (75.89), maa://37962 (90.24)"
$ws.Range("AF12").Value = "*maa://28932 (77.7), *maa://20106 (63.96), *maa://22769 (64.29)"
$ws.Range("D13").Value = "maa://24999 (92.13), maa://36673 (93.33), maa://25001 (85.71)"
$ws.Range("P13").Value = "maa://22676 (92.8), *maa://22583 (75.0), *maa://22500 (58.7)"
$ws.Range("D14").Value = "maa://30764 (89.09)"
$ws.Range("L14").Value = "maa://26245 (96.73), maa://21288 (96.3), maa://39841 (95.33), maa://36682 (97.44)"
$ws.Range("AF15").Value = "maa://21364 (80.78), *maa://36666 (77.48), *maa://22766 (68.64)"
$ws.Range("D16").Value = "maa://21441 (96.41), maa://36679 (94.44), maa://37650 (97.3)"
$ws.Range("D18").Value = "maa://24570 (97.36)"
$ws.Range("L18").Value = "maa://22466 (90.42), *maa://22732 (50.55)"
$ws.Range("T19").Value = "maa://24386 (99.15)"
$ws.Range("AF19").Value = "*maa://21663 (63.89)"
$ws.Range("L20").Value = "maa://41331 (85.16)"
$ws.Range("P23").Value = "maa://30587 (91.92), *maa://29748 (75.97), ***maa://29785 (16.18), *maa://37566 (76.92)"
$ws.Range("X24").Value = "maa://29988 (84.31), maa://23504 (93.19), **maa://22892 (40.14), *maa://25141 (77.1), *maa://36663 (76.92), ***maa://22815 (23.08)"
$ws.Range("H25").Value = "*maa://29063 (73.17), *maa://25311 (74.04), ***maa://22725 (4.84), *maa://45047 (62.5)"
$ws.Range("AB26").Value = "maa://42235 (94.17)"
$ws.Range("AF26").Value = "maa://30511 (80.95), *maa://29760 (60.0)"
$ws.Range("AF27").Value = "maa://24023 (97.3)"
$ws.Range("D28").Value = "maa://24465 (90.94), maa://25725 (83.91)"
$ws.Range("X28").Value = "maa://39929 (90.77), maa://41749 (90.7), ***maa://39723 (13.89)"
$ws.Range("AF28").Value = "maa://36660 (92.31), *maa://36701 (66.67)"
$ws.Range("L29").Value = "maa://28432 (93.49), *maa://28440 (80.0), maa://31400 (98.81), *maa://28650 (71.43)"
$ws.Range("AF29").Value = "*maa://24080 (68.85), maa://42865 (81.25), ***maa://34960 (8.33)"
$ws.Range("P30").Value = "maa://21442 (99.11)"
$ws.Range("AB30").Value = "maa://42979 (97.01), maa://45822 (100.0), *maa://45045 (80.0)"
$ws.Range("L31").Value = "maa://35926 (93.49), maa://36258 (85.47), *maa://43904 (72.73)"
$ws.Range("T32").Value = "maa://42859 (95.97), maa://41108 (88.0), maa://41238 (97.12), maa://45523 (100.0)"
$ws.Range("L35").Value = "maa://41296 (96.27)"
$ws.Range("L37").Value = "maa://45718 (97.89), *maa://47069 (69.23), maa://45789 (100.0)"
$ws.Range("H39").Value = "maa://36670 (89.11), maa://25199 (84.82), maa://30434 (91.46), ***maa://25036 (16.0), maa://45059 (83.33), *maa://44165 (66.67)"
$ws.Range("T39").Value = "maa://45788 (81.44), maa://47079 (94.12), *maa://45790 (75.0)"
$ws.Range("P40").Value = "maa://23278 (95.56), maa://21386 (95.77), maa://36664 (89.29), maa://45550 (100.0)"
$ws.Range("H45").Value = "maa://21229 (84.82), maa://30807 (95.65), *maa://22767 (55.0), ***maa://20796 (13.79), maa://42459 (85.71)"
$ws.Range("H53").Value = "maa://32534 (94.1), **maa://32434 (33.33)"
$ws.Range("H55").Value = "maa://32532 (91.91)"
$ws.Range("H58").Value = "*maa://37964 (61.54)"
